$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.7
$ws.Range("C3").Value = 9.3

# Update active selection to C3
$ws.Range("C3").Select()
